$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FBA79697","TD310","B0CQX4K9P5",44,59623.52),
    @("FBA77113","K1","B01ISNU3X4",35,53892.28),
    @("FBA79113","TC310","B0BTCXQQ6M",18,33544.08),
    @("FBA79696","TD310+","B0CQX3VB1R",15,22868.7),
    @("FBA79114","TC310+","B0CCV74CL7",7,17790.7),
    @("FBA79260","G11","B07GVGMW59",8,17281.35),
    @("FBA79111","TD510","B0BRKFP94K",4,13555.92),
    @("FBA79116","TC-777 PRO","B0BYHHSLPC",8,13552.53),
    @("FBA77117","S20","B078WNW4YW",6,13215.24),
    @("FBA77101","TC-777","B07WLWN2ZT",7,13044.92),
    @("FBA77106","T20","B082W4B7SX",4,8132.21),
    @("FBA77105","T30","B089FVQD3Z",3,7370.34),
    @("FBA77111","TC30","B08CVP2HXP",3,5844.93),
    @("FBA79574","TC30S","B0B4WTHLX5",2,4235.6),
    @("FBA77114","TC-2030","B07TSN2H9D",1,3643.22),
    @("FBA77110","TM20","B08NDB5NWP",0,0)
)

$fmt = $ws.Cells.Item(24, 20).NumberFormat

$row = 25
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 16).Value = $item[3]
    $ws.Cells.Item($row, 20).Value = $item[4]
    $ws.Cells.Item($row, 20).NumberFormat = $fmt
    $row++
}

[void]$ws.Range("A25").Select()
[void]($excel.ActiveWindow.ScrollRow = 2)
